# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (E2) and "Correspond Handback DateTime" (H2)
# timestamps on the zh-cn and de-de worksheets to reflect the new handback run times.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-18 16:57:03"
$wsZhCn.Range("H2").Value = "2016-03-18 16:57:16"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-18 16:57:07"
$wsDeDe.Range("H2").Value = "2016-03-18 16:57:22"
